$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '43.283.26'
$ws.Range("E2").Value = '  -1.27%  '
$ws.Range("D3").Value = '2.273.27'
$ws.Range("E3").Value = '  -1.87%  '
$ws.Range("E4").Value = '  -0.03%  '
$ws.Range("E5").Value = '  -2.09%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '264.80'
$ws.Range("E6").Value = '  -2.22%  '
$ws.Range("E7").Value = '  -1.17%  '
$ws.Range("E8").Value = '  +0.33%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.609'
$ws.Range("E9").Value = '  -2.41%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '47.92'
$ws.Range("E10").Value = '  +0.71%  '
$ws.Range("E11").Value = '  -1.35%  '
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '8.79'
$ws.Range("E12").Value = '  -1.17%  '
$ws.Range("E13").Value = '  +0.87%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '15.47'
$ws.Range("E14").Value = '  -1.83%  '
$ws.Range("D15").Value = '2.616.33'
$ws.Range("E15").Value = '  -1.77%  '
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '0.854'
$ws.Range("E16").Value = '  -1.03%  '
$ws.Range("D17").Value = '2.268.69'
$ws.Range("E17").Value = '  -2.25%  '
$ws.Range("D18").Value = '43.175.39'
$ws.Range("E18").Value = '  -1.48%  '
$ws.Range("E19").Value = '  -2.42%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '6.77'
$ws.Range("E20").Value = '  +1.23%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '71.21'
$ws.Range("E21").Value = '  -2.25%  '
$ws.Range("E22").Value = '  +0.74%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '231.80'
$ws.Range("E23").Value = '  -1.16%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '9.73'
$ws.Range("E24").Value = '  +2.56%  '
$ws.Range("E25").Value = '  -1.98%  '
$ws.Range("E26").Value = '  +0.34%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '11.31'
$ws.Range("E27").Value = '  -1.60%  '
$ws.Range("E28").Value = '  -1.03%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '40.35'
$ws.Range("E29").Value = '  -6.02%  '
$ws.Range("B30").Value = 'WEMIXToken'
$ws.Range("C30").Value = 'https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix'
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '3.32'
$ws.Range("E30").Value = '  -3.57%  '
$ws.Range("B31").Value = 'Toncoin'
$ws.Range("C31").Value = 'https://coinranking.com/coin/67YlI0K1b+toncoin-ton'
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '2.25'
$ws.Range("E31").Value = '  -1.55%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '171.86'
$ws.Range("E32").Value = '  -3.37%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '21.32'
$ws.Range("E33").Value = '  -3.06%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '0.0905'
$ws.Range("E34").Value = '  -2.98%  '
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '5.74'
$ws.Range("E35").Value = '  +2.56%  '
$ws.Range("E36").Value = '  +0.33%  '
$ws.Range("E37").Value = '  -1.84%  '
$ws.Range("E38").Value = '  -1.36%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '3.84'
$ws.Range("E39").Value = '  -3.64%  '
$ws.Range("E40").Value = '  -6.12%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '2.61'
$ws.Range("E41").Value = '  +8.39%  '
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '77.28'
$ws.Range("E42").Value = '  +9.61%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '14.02'
$ws.Range("E43").Value = '  +10.52%  '
$ws.Range("E44").Value = '  -3.77%  '
$ws.Range("E45").Value = '  +2.18%  '
$ws.Range("E46").Value = '  +0.18%  '
$ws.Range("E47").Value = '  -2.42%  '
$ws.Range("E48").Value = '  -2.47%  '
$ws.Range("E49").Value = '  -1.55%  '
$ws.Range("B50").Value = 'TrustWalletToken'
$ws.Range("C50").Value = 'https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt'
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '1.25'
$ws.Range("E50").Value = '  +1.34%  '
$ws.Range("B51").Value = 'Aave'
$ws.Range("C51").Value = 'https://coinranking.com/coin/ixgUfzmLR+aave-aave'
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '101.01'
$ws.Range("E51").Value = '  +0.49%  '
